$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("partidas")

# Row 13 corresponds to "Partida" 12: fill in the result of the match
# E13 (RF - score string), G13/H13 (goals), I13 (Status)
$ws.Range("E13").Value = "7x3"
$ws.Range("G13").Value = 7
$ws.Range("H13").Value = 3
$ws.Range("I13").Value = "Finalizado"

# Update the active selection to J13
$ws.Range("J13").Select()
